# Update "想去人数" (interested count) values for three events that
# appear both on the "展览" sheet and the "全部类型" sheet.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F3").Value = 1665
    $ws.Range("F5").Value = 6220
    $ws.Range("F6").Value = 53
}
